# env 7010 update test objects to follow global variable

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("login")

# Add new PTJ header column (C1), matching the style/value pattern used for A1/B1
$ws.Range("C1").Value = "PTJ"

# Populate rows 5-8 with new test data (IC number in col A, fixed password-like value in col B)
$data = @(
    @(810213016116, "8SQVv/p9jVScEs4/2CZsLw=="),
    @(600803015724, "8SQVv/p9jVScEs4/2CZsLw=="),
    @(820806015126, "8SQVv/p9jVScEs4/2CZsLw=="),
    @(780613015177, "8SQVv/p9jVScEs4/2CZsLw==")
)

$r = 5
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r++
}

# Adjust column A width to fit the new, wider IC number header/content
$ws.Columns.Item(1).ColumnWidth = 19

# Update the selected cell to reflect the new working position
$ws.Range("B8").Select()
